# Fixing missing color in slide
#
# Three "TextBox 34" / "TextBox 59" legend textboxes (on slides 3, 4 and 5)
# list the colors used in the DCF diagrams. The legend text for the first
# bullet ("Blue: ...") was already updated to "Blue: Original DCF Transmit"
# and a new "Yellow: Original DCF Backoff" bullet was already added on
# slide 3, but that same update was missed on slides 4 and 5 - and the
# legend box on all three slides needs to move up a bit to make room.
#
# This script:
#   1. Nudges the legend textbox up on slides 3, 4 and 5 (same new Top).
#   2. Adds the missing "Blue: ... Transmit" / "Yellow: ... Backoff" text
#      (and the corresponding box-height growth) to the legend boxes on
#      slides 4 and 5 so they match the one already fixed on slide 3.

$p = $ppt.ActivePresentation

# Point value that serializes to exactly 1934775 EMU. Shape.Top/.Height
# are single-precision floats, so the value is chosen to round-trip to
# the exact target EMU through that float conversion.
$newTop = 152.3445

function Set-LegendBoxTop($slideIndex, $shapeIndex) {
    $shape = $p.Slides.Item($slideIndex).Shapes.Item($shapeIndex)
    $shape.Top = $newTop
}

function Update-LegendText($slideIndex, $shapeIndex, $newHeight) {
    $shape = $p.Slides.Item($slideIndex).Shapes.Item($shapeIndex)
    $textRange = $shape.TextFrame.TextRange

    # "Blue: Original DCF" -> "Blue: Original DCF Transmit". Editing the
    # run's Text (rather than the paragraph's) keeps it as a single run.
    $firstRun = $textRange.Paragraphs(1).Runs(1)
    $firstRun.Text = "Blue: Original DCF Transmit"

    # Insert the new "Yellow: Original DCF Backoff" bullet right after it
    # as its own paragraph (it picks up the same bullet pPr as its
    # neighbors automatically).
    $firstRun.InsertAfter("`rYellow: Original DCF Backoff") | Out-Null

    # Split "Backoff" into its own run, mirroring the source formatting:
    # two runs, "Yellow: Original DCF " and "Backoff".
    $newPara = $textRange.Paragraphs(2)
    $backoffStart = $newPara.Text.IndexOf("Backoff") + 1
    $backoffRun = $newPara.Characters($backoffStart, 7)
    $backoffRun.Text = "Backoff"

    # The box has spAutoFit, so adding the paragraph above already grew
    # it - set Top/Height explicitly afterwards so the final values are
    # exact (not just whatever autofit happened to compute).
    $shape.Top = $newTop
    $shape.Height = $newHeight
}

# Slide 3 ("Compressible DCF"): legend box already has the "Transmit" /
# "Backoff" text - it only needs to move up.
Set-LegendBoxTop 3 32

# Slide 4 ("Adding Interarrival"): legend box needs the text fix and to
# grow taller to fit the extra bullet line.
Update-LegendText 4 2 104.2078

# Slide 5 ("Adding Interarrival" / packet size): same fix; the box ends
# up one line taller than slide 4's (it already had a 4th bullet before
# this edit, vs. slide 4's 3).
Update-LegendText 5 7 128.44221
